# Restore cell C10 on the "Rules" sheet from 18 to 1 (re-saved workbook value change).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
